$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by the Overview sheet (B/C columns) and the
#    per-locale sheets' "Status" (C) column for rows 2 and 3.
# ---------------------------------------------------------------------------
$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

$locales = @("zh-cn", "de-de")
foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc)
    $ws.Range("C2").Value = $statusNew
    $ws.Range("C3").Value = $statusNew
}

# ---------------------------------------------------------------------------
# 2. For zh-cn and de-de sheets: populate new "Latest Target File" (F) and
#    "Latest Handback File" (G) columns for rows 2 and 3, with hyperlinks
#    mirroring the existing Source File (A) / Latest Handoff File (D) links.
#    Also set the "Latest Handback DateTime" (H) which used to be the
#    zero-date placeholder.
# ---------------------------------------------------------------------------

# zh-cn sheet
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("F2").Value = "7796e927-103a-4fc2-bf15-49a58ff275e9.md"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c8ee7adcf8cb09fd56764a985b6ce524e7fe416/e2e/7796e927-103a-4fc2-bf15-49a58ff275e9.md", [Type]::Missing, [Type]::Missing, "7796e927-103a-4fc2-bf15-49a58ff275e9.md")

$ws.Range("G2").Value = "7796e927-103a-4fc2-bf15-49a58ff275e9.1abda1edfffe9f485eeb9850be2a73db422b8859.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf3e1c3675cfd530237324c6d8f3d59e7b754b15/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7796e927-103a-4fc2-bf15-49a58ff275e9.1abda1edfffe9f485eeb9850be2a73db422b8859.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "7796e927-103a-4fc2-bf15-49a58ff275e9.1abda1edfffe9f485eeb9850be2a73db422b8859.zh-cn.xlf")

$ws.Range("H2").Value = "2016-03-14 01:18:13"

$ws.Range("F3").Value = "8816f159-8b35-41a5-8cc2-8925efe070e2.md"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c8ee7adcf8cb09fd56764a985b6ce524e7fe416/e2e/8816f159-8b35-41a5-8cc2-8925efe070e2.md", [Type]::Missing, [Type]::Missing, "8816f159-8b35-41a5-8cc2-8925efe070e2.md")

$ws.Range("G3").Value = "8816f159-8b35-41a5-8cc2-8925efe070e2.cc02c60300447251f360b678bc7461a40bb5c447.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf3e1c3675cfd530237324c6d8f3d59e7b754b15/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8816f159-8b35-41a5-8cc2-8925efe070e2.cc02c60300447251f360b678bc7461a40bb5c447.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "8816f159-8b35-41a5-8cc2-8925efe070e2.cc02c60300447251f360b678bc7461a40bb5c447.zh-cn.xlf")

$ws.Range("H3").Value = "2016-03-14 01:18:13"

foreach ($ref in @("F2", "G2", "F3", "G3")) {
    $ws.Range($ref).Style = "HyperLink"
}

# de-de sheet
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("F2").Value = "7796e927-103a-4fc2-bf15-49a58ff275e9.md"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c8ee7adcf8cb09fd56764a985b6ce524e7fe416/e2e/7796e927-103a-4fc2-bf15-49a58ff275e9.md", [Type]::Missing, [Type]::Missing, "7796e927-103a-4fc2-bf15-49a58ff275e9.md")

$ws.Range("G2").Value = "7796e927-103a-4fc2-bf15-49a58ff275e9.1abda1edfffe9f485eeb9850be2a73db422b8859.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/effffe833a274f96efbfe5f35f48f9418396a2b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7796e927-103a-4fc2-bf15-49a58ff275e9.1abda1edfffe9f485eeb9850be2a73db422b8859.de-de.xlf", [Type]::Missing, [Type]::Missing, "7796e927-103a-4fc2-bf15-49a58ff275e9.1abda1edfffe9f485eeb9850be2a73db422b8859.de-de.xlf")

$ws.Range("H2").Value = "2016-03-14 01:18:19"

$ws.Range("F3").Value = "8816f159-8b35-41a5-8cc2-8925efe070e2.md"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c8ee7adcf8cb09fd56764a985b6ce524e7fe416/e2e/8816f159-8b35-41a5-8cc2-8925efe070e2.md", [Type]::Missing, [Type]::Missing, "8816f159-8b35-41a5-8cc2-8925efe070e2.md")

$ws.Range("G3").Value = "8816f159-8b35-41a5-8cc2-8925efe070e2.cc02c60300447251f360b678bc7461a40bb5c447.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/effffe833a274f96efbfe5f35f48f9418396a2b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8816f159-8b35-41a5-8cc2-8925efe070e2.cc02c60300447251f360b678bc7461a40bb5c447.de-de.xlf", [Type]::Missing, [Type]::Missing, "8816f159-8b35-41a5-8cc2-8925efe070e2.cc02c60300447251f360b678bc7461a40bb5c447.de-de.xlf")

$ws.Range("H3").Value = "2016-03-14 01:18:19"

foreach ($ref in @("F2", "G2", "F3", "G3")) {
    $ws.Range($ref).Style = "HyperLink"
}
